$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the header row: "_old" columns become "_FV2210", "_new" columns become "_FV2304"
for ($c = 1; $c -le 21; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $val = $cell.Value()
    if ($val -ne $null) {
        $newVal = $val -replace "_old$", "_FV2210"
        $newVal = $newVal -replace "_new$", "_FV2304"
        if ($newVal -ne $val) {
            $cell.Value = $newVal
        }
    }
}

# Turn the used range into an Excel Table ("Table1") with headers
$tableRange = $ws.Range("A1:U55")
$tbl = $ws.ListObjects.Add(1, $tableRange, [System.Reflection.Missing]::Value, 1)
$tbl.Name = "Table1"
$tbl.TableStyle = ""

# Freeze the header row (top row)
$ws.Activate() | Out-Null
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
